# Add a new "2025-06-23" attendance column (U) and update the
# Total / Attendance % summary columns to account for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell U1: "2025-06-23" (styled like the other header cells) ---
$ws.Range("U1").Value = "'2025-06-23"
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2 (Abhishek Pathak): absent on the new date, totals updated ---
$ws.Range("S2").Value = 16
$ws.Range("T2").Value = 6.2
$ws.Range("U2").Value = "❌"

# --- Row 3 (Shubham Pitekar): absent on the new date, totals updated ---
$ws.Range("S3").Value = 16
$ws.Range("U3").Value = "❌"
